$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits -----------------------------------------------------

# Row 4: false negative corrected from 1 to 0
$ws.Range("D4").Value = 0

# Row 5 / Row 6: remove the old per-row notes (these notes no longer apply
# with the smaller uploaded data set)
$ws.Range("E5").ClearContents()
$ws.Range("E6").ClearContents()

# Row 8: fill in counts + note now that more data is available
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = "CP count based on looking at obj image for this point and below"

# Row 9
$ws.Range("B9").Value = 9
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0

# Row 10
$ws.Range("B10").Value = 13
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0

# Row 11
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0

# Row 12
$ws.Range("B12").Value = 16
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = "1 possible false pos (out of focus) lower right"

# Row 13
$ws.Range("B13").Value = 12
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = "Edge objects inconsistent -- here a cell was detected that is at edge (partial nucleus) but in #43 a cell barely on edge was eliminated"

# Row 14
$ws.Range("B14").Value = 15
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0

# Row 15
$ws.Range("B15").Value = 10
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0

# Row 16
$ws.Range("B16").Value = 14
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = "false pos = elongated junk"

# Row 17
$ws.Range("B17").Value = 9
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 0

# --- View state -------------------------------------------------------
# Scroll the frozen sheet down and leave the selection on E7, matching
# where the author was working when the data above was entered.
$win = $excel.ActiveWindow
$win.ScrollRow = 9
$win.ScrollColumn = 1
$ws.Range("E7").Select()
